$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the missing Tuesday/Wednesday hours for the week of row 5
$ws.Range("E5").Value = 3.5
$ws.Range("F5").Value = 3

# Update the active cell / selection to G5 (matches the recorded view state)
$ws.Range("G5").Select()
